# Atualizado por script em 05-11-2023 14:45
# Updates betting-odds rows for Liga Portugal 2 (2023-2024):
#   - row 74 <-> row 75 swap (FC Porto B-Feirense / Mafra-Leixoes)
#   - row 77 <-> row 78 swap (Benfica B-Penafiel / Nacional-Santa Clara)
#   - new row 79 appended (AVS - Oliveirense)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 74: becomes Mafra vs Leixoes ----
$ws.Range("F74").Value = "Mafra"
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = "Leixoes"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1.88
$ws.Range("K74").Value = "01/11/2023 16:12"
$ws.Range("L74").Value = 1.93
$ws.Range("M74").Value = "04/11/2023 11:48"
$ws.Range("N74").Value = 3.73
$ws.Range("O74").Value = "01/11/2023 16:12"
$ws.Range("P74").Value = 3.55
$ws.Range("Q74").Value = "04/11/2023 11:51"
$ws.Range("R74").Value = 3.8
$ws.Range("S74").Value = "01/11/2023 16:12"
$ws.Range("T74").Value = 4.21
$ws.Range("U74").Value = "04/11/2023 11:51"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/mafra-leixoes/YiBBPnTT/"

# ---- Row 75: becomes FC Porto B vs Feirense ----
$ws.Range("F75").Value = "FC Porto B"
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = "Feirense"
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 1.98
$ws.Range("K75").Value = "01/11/2023 16:12"
$ws.Range("L75").Value = 1.84
$ws.Range("M75").Value = "04/11/2023 11:59"
$ws.Range("N75").Value = 3.57
$ws.Range("O75").Value = "01/11/2023 16:12"
$ws.Range("P75").Value = 3.78
$ws.Range("Q75").Value = "04/11/2023 11:59"
$ws.Range("R75").Value = 3.87
$ws.Range("S75").Value = "01/11/2023 16:12"
$ws.Range("T75").Value = 4.41
$ws.Range("U75").Value = "04/11/2023 11:58"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/fc-porto-feirense/jTL6QSDN/"

# ---- Row 77: becomes Nacional vs Santa Clara ----
$ws.Range("F77").Value = "Nacional"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Santa Clara"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 2.98
$ws.Range("K77").Value = "01/11/2023 16:12"
$ws.Range("L77").Value = 2.81
$ws.Range("M77").Value = "04/11/2023 18:58"
$ws.Range("N77").Value = 3.27
$ws.Range("O77").Value = "01/11/2023 16:12"
$ws.Range("P77").Value = 3.23
$ws.Range("Q77").Value = "04/11/2023 18:52"
$ws.Range("R77").Value = 2.39
$ws.Range("S77").Value = "01/11/2023 16:12"
$ws.Range("T77").Value = 2.72
$ws.Range("U77").Value = "04/11/2023 18:52"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/"

# ---- Row 78: becomes Benfica B vs Penafiel ----
$ws.Range("F78").Value = "Benfica B"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Penafiel"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.91
$ws.Range("K78").Value = "29/10/2023 16:42"
$ws.Range("L78").Value = 2.26
$ws.Range("M78").Value = "04/11/2023 18:53"
$ws.Range("N78").Value = 3.67
$ws.Range("O78").Value = "29/10/2023 16:42"
$ws.Range("P78").Value = 3.59
$ws.Range("Q78").Value = "04/11/2023 18:53"
$ws.Range("R78").Value = 4.01
$ws.Range("S78").Value = "29/10/2023 16:42"
$ws.Range("T78").Value = 3.2
$ws.Range("U78").Value = "04/11/2023 18:53"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/"

# ---- Row 79: new row, AVS vs Oliveirense ----
# Copy formatting from row 78 (the current last data row) so the new row
# matches the sheet's existing style (bold/bordered index column, date
# number-format column, etc.) before filling in values.
$ws.Range("A78:V78").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "portugal"
$ws.Range("C79").Value = "liga-portugal-2"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45235.5
$ws.Range("F79").Value = "AVS"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = "Oliveirense"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 1.51
$ws.Range("K79").Value = "01/11/2023 16:12"
$ws.Range("L79").Value = 1.51
$ws.Range("M79").Value = "05/11/2023 11:57"
$ws.Range("N79").Value = 4.39
$ws.Range("O79").Value = "01/11/2023 16:12"
$ws.Range("P79").Value = 4.49
$ws.Range("Q79").Value = "05/11/2023 11:57"
$ws.Range("R79").Value = 6.33
$ws.Range("S79").Value = "01/11/2023 16:12"
$ws.Range("T79").Value = 6.55
$ws.Range("U79").Value = "05/11/2023 11:57"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/avs-oliveirense/4xkjWWqn/"
